$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = "payment_plan_days_mov_avg_m3"
$ws.Cells.Item(4, 2).Value = 30
$ws.Cells.Item(5, 1).Value = "payment_plan_days_mov_avg_m6"
$ws.Cells.Item(5, 2).Value = 30
$ws.Cells.Item(6, 1).Value = "payment_plan_days_mov_max_m3"
$ws.Cells.Item(6, 2).Value = 30
$ws.Cells.Item(7, 1).Value = "payment_plan_days_mov_max_m6"
$ws.Cells.Item(7, 2).Value = 30
$ws.Cells.Item(8, 1).Value = "payment_plan_days_mov_min_m3"
$ws.Cells.Item(8, 2).Value = 30
$ws.Cells.Item(9, 1).Value = "payment_plan_days_mov_min_m6"
$ws.Cells.Item(9, 2).Value = 30
$ws.Cells.Item(10, 1).Value = "actual_amount_paid"
$ws.Cells.Item(10, 2).Value = 149
$ws.Cells.Item(11, 1).Value = "actual_amount_paid_mov_avg_m3"
$ws.Cells.Item(11, 2).Value = 149
$ws.Cells.Item(12, 1).Value = "actual_amount_paid_mov_avg_m6"
$ws.Cells.Item(12, 2).Value = 149
$ws.Cells.Item(13, 1).Value = "actual_amount_paid_mov_max_m3"
$ws.Cells.Item(13, 2).Value = 149
$ws.Cells.Item(14, 1).Value = "actual_amount_paid_mov_max_m6"
$ws.Cells.Item(14, 2).Value = 149
$ws.Cells.Item(15, 1).Value = "actual_amount_paid_mov_min_m3"
$ws.Cells.Item(15, 2).Value = 149
$ws.Cells.Item(16, 1).Value = "actual_amount_paid_mov_min_m6"
$ws.Cells.Item(16, 2).Value = 149
$ws.Cells.Item(17, 1).Value = "num_25"
$ws.Cells.Item(17, 2).Value = 57
$ws.Cells.Item(18, 1).Value = "num_25_mov_avg_m3"
$ws.Cells.Item(18, 2).Value = 63
$ws.Cells.Item(19, 1).Value = "num_25_mov_avg_m6"
$ws.Cells.Item(19, 2).Value = 66.5
$ws.Cells.Item(20, 1).Value = "num_25_mov_max_m3"
$ws.Cells.Item(20, 2).Value = 93
$ws.Cells.Item(21, 1).Value = "num_25_mov_max_m6"
$ws.Cells.Item(21, 2).Value = 122
$ws.Cells.Item(22, 1).Value = "num_25_mov_min_m3"
$ws.Cells.Item(22, 2).Value = 63
$ws.Cells.Item(23, 1).Value = "num_25_mov_min_m6"
$ws.Cells.Item(23, 2).Value = 66.5
$ws.Cells.Item(24, 1).Value = "num_50"
$ws.Cells.Item(24, 2).Value = 15
$ws.Cells.Item(25, 1).Value = "num_50_mov_avg_m3"
$ws.Cells.Item(25, 2).Value = 17
$ws.Cells.Item(26, 1).Value = "num_50_mov_avg_m6"
$ws.Cells.Item(26, 2).Value = 18
$ws.Cells.Item(27, 1).Value = "num_50_mov_max_m3"
$ws.Cells.Item(27, 2).Value = 25
$ws.Cells.Item(28, 1).Value = "num_50_mov_max_m6"
$ws.Cells.Item(28, 2).Value = 33
$ws.Cells.Item(29, 1).Value = "num_50_mov_min_m3"
$ws.Cells.Item(29, 2).Value = 17
$ws.Cells.Item(30, 1).Value = "num_50_mov_min_m6"
$ws.Cells.Item(30, 2).Value = 18
$ws.Cells.Item(31, 1).Value = "num_75"
$ws.Cells.Item(31, 2).Value = 10
$ws.Cells.Item(32, 1).Value = "num_75_mov_avg_m3"
$ws.Cells.Item(32, 2).Value = 11
$ws.Cells.Item(33, 1).Value = "num_75_mov_avg_m6"
$ws.Cells.Item(33, 2).Value = 11.5
$ws.Cells.Item(34, 1).Value = "num_75_mov_max_m3"
$ws.Cells.Item(34, 2).Value = 16
$ws.Cells.Item(35, 1).Value = "num_75_mov_max_m6"
$ws.Cells.Item(35, 2).Value = 21
$ws.Cells.Item(36, 1).Value = "num_75_mov_min_m3"
$ws.Cells.Item(36, 2).Value = 11
$ws.Cells.Item(37, 1).Value = "num_75_mov_min_m6"
$ws.Cells.Item(37, 2).Value = 11.5
$ws.Cells.Item(38, 1).Value = "num_985"
$ws.Cells.Item(38, 2).Value = 11
$ws.Cells.Item(39, 1).Value = "num_985_mov_avg_m3"
$ws.Cells.Item(39, 2).Value = 11.33333333333333
$ws.Cells.Item(40, 1).Value = "num_985_mov_avg_m6"
$ws.Cells.Item(40, 2).Value = 11.66666666666667
$ws.Cells.Item(41, 1).Value = "num_985_mov_max_m3"
$ws.Cells.Item(41, 2).Value = 17
$ws.Cells.Item(42, 1).Value = "num_985_mov_max_m6"
$ws.Cells.Item(42, 2).Value = 21
$ws.Cells.Item(43, 1).Value = "num_985_mov_min_m3"
$ws.Cells.Item(43, 2).Value = 11.33333333333333
$ws.Cells.Item(44, 1).Value = "num_985_mov_min_m6"
$ws.Cells.Item(44, 2).Value = 11.66666666666667
$ws.Cells.Item(45, 1).Value = "num_100"
$ws.Cells.Item(45, 2).Value = 298
$ws.Cells.Item(46, 1).Value = "num_100_mov_avg_m3"
$ws.Cells.Item(46, 2).Value = 315.3333333333333
$ws.Cells.Item(47, 1).Value = "num_100_mov_avg_m6"
$ws.Cells.Item(47, 2).Value = 327.1666666666667
$ws.Cells.Item(48, 1).Value = "num_100_mov_max_m3"
$ws.Cells.Item(48, 2).Value = 441
$ws.Cells.Item(49, 1).Value = "num_100_mov_max_m6"
$ws.Cells.Item(49, 2).Value = 548
$ws.Cells.Item(50, 1).Value = "num_100_mov_min_m3"
$ws.Cells.Item(50, 2).Value = 315.3333333333333
$ws.Cells.Item(51, 1).Value = "num_100_mov_min_m6"
$ws.Cells.Item(51, 2).Value = 327.1666666666667
$ws.Cells.Item(52, 1).Value = "num_unq"
$ws.Cells.Item(52, 2).Value = 334
$ws.Cells.Item(53, 1).Value = "num_unq_mov_avg_m3"
$ws.Cells.Item(53, 2).Value = 350
$ws.Cells.Item(54, 1).Value = "num_unq_mov_avg_m6"
$ws.Cells.Item(54, 2).Value = 361.8333333333333
$ws.Cells.Item(55, 1).Value = "num_unq_mov_max_m3"
$ws.Cells.Item(55, 2).Value = 483
$ws.Cells.Item(56, 1).Value = "num_unq_mov_max_m6"
$ws.Cells.Item(56, 2).Value = 593
$ws.Cells.Item(57, 1).Value = "num_unq_mov_min_m3"
$ws.Cells.Item(57, 2).Value = 350
$ws.Cells.Item(58, 1).Value = "num_unq_mov_min_m6"
$ws.Cells.Item(58, 2).Value = 361.8333333333333
$ws.Cells.Item(59, 1).Value = "total_secs"
$ws.Cells.Item(59, 2).Value = 82250.04399999999
$ws.Cells.Item(60, 1).Value = "%num_more_than_50"
$ws.Cells.Item(60, 2).Value = 79.40000000000001
$ws.Cells.Item(61, 1).Value = "%num_more_than_50_mov_avg_m3"
$ws.Cells.Item(61, 2).Value = 78.13000000000001
$ws.Cells.Item(62, 1).Value = "%num_more_than_50_mov_avg_m6"
$ws.Cells.Item(62, 2).Value = 77.45999999999999
$ws.Cells.Item(63, 1).Value = "%num_more_than_50_mov_max_m3"
$ws.Cells.Item(63, 2).Value = 85.59999999999999
$ws.Cells.Item(64, 1).Value = "%num_more_than_50_mov_max_m6"
$ws.Cells.Item(64, 2).Value = 88.37
$ws.Cells.Item(65, 1).Value = "%num_more_than_50_mov_min_m3"
$ws.Cells.Item(65, 2).Value = 78.13000000000001
$ws.Cells.Item(66, 1).Value = "%num_more_than_50_mov_min_m6"
$ws.Cells.Item(66, 2).Value = 77.45999999999999
$ws.Cells.Item(67, 1).Value = "months_as_a_registered"
$ws.Cells.Item(67, 2).Value = 30